$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("B1").Value = "Activation Date(dd/mm/yyyy)"
$ws.Range("B1").Select()
